# "article 93 is live"
#
# Row 7 of Sheet1 holds one "card" per column; the card in column I was the
# blog card pointing at article 90. That article has now been published as
# "ser: 93", so the card content is updated in place (Excel's shared-string
# table will naturally drop the now-unused "ser: 90" entry and append the
# new "ser: 93" entry on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 93"

# The editor's cursor ends up on C7 after making the change.
$ws.Range("C7").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 2
